# Data118Bus.xlsx - "Test 118 bus system"
# Update Apparatus-sheet load/impedance figures, drop a few helper formulas
# that became unnecessary once their values were hard-coded, flip the
# "Enable (create simulink model)" switch off on the Advance sheet, and
# leave the selections where the author last left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Apparatus sheet: update load / impedance values
# ---------------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("Apparatus")

# D column: 100 -> 200 (bus-load rows)
$rowsD100to200 = @(3,6,10,17,20,26,27,33,36,38,42,48,51,56,57,61,63,67,71,74,75,78,82)
foreach ($r in $rowsD100to200) { $wsApp.Cells.Item($r, 4).Value = 200 }

# D column: 50 -> 200 (bus-load rows, second block)
$rowsD50to200 = @(89,91,92,101,102,105,106,109,112,113)
foreach ($r in $rowsD50to200) { $wsApp.Cells.Item($r, 4).Value = 200 }

# D column: 2E-3 -> 5E-3 (apparatus sub-rows)
$rowsD2e3to5e3 = @(8,12,14,21,28,29,34,44,58,64,68,72,76,79,87,93,94,107,114,115,118)
foreach ($r in $rowsD2e3to5e3) { $wsApp.Cells.Item($r, 4).Value = 0.005 }

# F column: these sub-rows no longer need their "=E/10" helper formula -
# hard-code the literal result instead (1E-3)
foreach ($r in $rowsD2e3to5e3) { $wsApp.Cells.Item($r, 6).Value = 0.001 }

# F6 also loses its "=E6/5" helper formula, keeping the literal result (1E-3)
$wsApp.Cells.Item(6, 6).Value = 0.001

# G/H columns: 10 -> 5 for the sub-rows whose F formula was just dropped
# (rows 87,93,94,107,114,115,118 keep G/H at 10 - unchanged)
$rowsGH10to5 = @(8,12,14,21,28,29,34,44,58,64,68,72,76,79)
foreach ($r in $rowsGH10to5) {
    $wsApp.Cells.Item($r, 7).Value = 5
    $wsApp.Cells.Item($r, 8).Value = 5
}

# ---------------------------------------------------------------------
# Advance sheet: turn off "Enable (create simulink model)"
# ---------------------------------------------------------------------
$wsAdv = $wb.Worksheets.Item("Advance")
$wsAdv.Range("B8").Value = 0

# ---------------------------------------------------------------------
# Restore cursor positions the author left behind. NetworkLine_IEEE is
# touched first so the final Activate()/Select() below leaves Apparatus
# as the selected tab (matches workbook.xml's activeTab).
# ---------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item("NetworkLine_IEEE")
$wsLine.Activate()
$wsLine.Range("A8").Select()

$wsApp.Activate()
$wsApp.Range("F7").Select()
